# Ch2 Protocol Suite and Key Agreement
#
# The title placeholders on the two "PSK" slides (the "no forward
# secrecy" and "full forward secrecy" key-agreement diagrams) previously
# had no explicit <p:spPr> transform and simply inherited their
# position/size from the slide layout/master. Give them an explicit
# position & size instead.
#
# PowerPoint's COM object model reports/accepts Shape.Left/Top/Width/
# Height in points (1 pt = 12700 EMU) and stores them as single
# precision floats, while the underlying OOXML stores whole EMU. Naive
# EMU/12700.0 division can therefore truncate to one EMU below the
# intended value once it round-trips through a float. EmuToPt nudges
# the point value up in tiny steps until converting it back through a
# (single precision) float and truncating reproduces the exact target
# EMU value, so the saved XML matches exactly.

function EmuToPt([double]$emu) {
    $pt = $emu / 12700.0
    for ($i = 0; $i -lt 2000; $i++) {
        $asSingle = [float]$pt
        $roundTripEmu = [math]::Floor([double]$asSingle * 12700.0)
        if ($roundTripEmu -eq $emu) {
            return $pt
        }
        $pt = $pt + 0.0000001
    }
    return $pt
}

$p = $ppt.ActivePresentation

# Slide 18 ("PSK: ...no forward secrecy"): title -> off x=795885 y=126800, ext cx=10515600 cy=1325563
$slide18 = $p.Slides.Item(18)
$title18 = $slide18.Shapes.Item(2)
$title18.Left   = EmuToPt 795885
$title18.Top    = EmuToPt 126800
$title18.Width  = EmuToPt 10515600
$title18.Height = EmuToPt 1325563

# Slide 19 ("PSK: ...full forward secrecy"): title -> off x=838200 y=46959, ext cx=10515600 cy=1325563
$slide19 = $p.Slides.Item(19)
$title19 = $slide19.Shapes.Item(5)
$title19.Left   = EmuToPt 838200
$title19.Top    = EmuToPt 46959
$title19.Width  = EmuToPt 10515600
$title19.Height = EmuToPt 1325563
